$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.655.64'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.859.62'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D5').Value = '''522.68'
$ws.Range('E5').Value = '  +5.95%  '
$ws.Range('D6').Value = '''141.06'
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('E7').Value = '  -2.33%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('E10').Value = '  -5.78%  '
$ws.Range('E11').Value = '  -8.60%  '
$ws.Range('D12').Value = '''41.68'
$ws.Range('E12').Value = '  -3.86%  '
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '4.474.36'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').Value = '''21.49'
$ws.Range('E15').Value = '  +8.12%  '
$ws.Range('D16').Value = '3.882.67'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '''14.07'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').Value = '68.642.42'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').Value = '''416.78'
$ws.Range('E21').Value = '  -5.79%  '
$ws.Range('D22').Value = '''3.51'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '''14.05'
$ws.Range('E23').Value = '  -3.22%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').Value = '''11.99'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''86.72'
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '''4.02'
$ws.Range('E26').Value = '  +6.04%  '
$ws.Range('E27').Value = '  -5.93%  '
$ws.Range('D28').Value = '''35.47'
$ws.Range('E28').Value = '  -4.56%  '
$ws.Range('D29').Value = '''13.41'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '''675.65'
$ws.Range('E30').Value = '  -4.35%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''6.97'
$ws.Range('E31').Value = '  +14.59%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.125'
$ws.Range('E32').Value = '  -5.13%  '
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('D34').Value = '''67.05'
$ws.Range('E34').Value = '  +8.85%  '
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('D36').Value = '0.0₃0851'
$ws.Range('E36').Value = '  -7.85%  '
$ws.Range('D37').Value = '''39.46'
$ws.Range('E37').Value = '  -3.61%  '
$ws.Range('E38').Value = '  +14.72%  '
$ws.Range('D39').Value = '''0.147'
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  -3.68%  '
$ws.Range('D43').Value = '''2.87'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('E44').Value = '  +5.03%  '
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '''0.141'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('B47').Value = 'FLOKI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D47').Value = '''0.000286'
$ws.Range('E47').Value = '  +17.19%  '
$ws.Range('D48').Value = '''3.00'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('E49').Value = '  -3.47%  '
$ws.Range('D50').Value = '''8.79'
$ws.Range('E50').Value = '  +4.00%  '
$ws.Range('D51').Value = '''142.63'
$ws.Range('E51').Value = '  -0.96%  '
